# Update enrollment figures in the "Inscricoes" sheet to match the new
# commit: increment several "Inscritos" (E), "Pagos" (F) and
# "Inscrições homologadas" (H) counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("F2").Value = 22
$ws.Range("H2").Value = 22

$ws.Range("F15").Value = 61
$ws.Range("H15").Value = 61

$ws.Range("E17").Value = 79
$ws.Range("F17").Value = 36
$ws.Range("H17").Value = 36

$ws.Range("E18").Value = 76

$ws.Range("E30").Value = 3

$ws.Range("F33").Value = 7
$ws.Range("H33").Value = 7

$ws.Range("E36").Value = 62

$ws.Range("E37").Value = 30
$ws.Range("F37").Value = 16
$ws.Range("H37").Value = 16

$ws.Range("E43").Value = 15

$ws.Range("E48").Value = 18

$ws.Range("F60").Value = 6
$ws.Range("H60").Value = 6

$ws.Range("E62").Value = 26

$ws.Range("E63").Value = 17

$ws.Range("E64").Value = 24

$ws.Range("E67").Value = 28

$ws.Range("F70").Value = 11
$ws.Range("H70").Value = 11

$ws.Range("E73").Value = 18

$ws.Range("E74").Value = 13

$ws.Range("E76").Value = 30

$ws.Range("E78").Value = 18
$ws.Range("F78").Value = 7
$ws.Range("H78").Value = 7

$ws.Range("F79").Value = 9
$ws.Range("H79").Value = 9

$ws.Range("F88").Value = 8
$ws.Range("H88").Value = 8
